$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B, C, D, E columns for rows 2-4 with new computed values
$ws.Range("B2").Value = -0.5604861728
$ws.Range("C2").Value = -224.73737164
$ws.Range("D2").Value = -225.29785782
$ws.Range("E2").Value = -224.5127374802

$ws.Range("B3").Value = -0.5691529588000001
$ws.Range("C3").Value = -224.66159912
$ws.Range("D3").Value = -225.23075208
$ws.Range("E3").Value = -224.5127374802

$ws.Range("B4").Value = -0.5726482441
$ws.Range("C4").Value = -224.64660074
$ws.Range("D4").Value = -225.21924899
$ws.Range("E4").Value = -224.5127374802
